# Update self-assessment ("Tự đánh giá", column F) scores for the
# category summary rows, and move the frozen-pane viewport / active
# selection down to the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Value = 2.75
$ws.Range("F21").Value = 0.25
$ws.Range("F25").Value = 1
$ws.Range("F28").Value = 0.25
$ws.Range("F36").Value = 0.25
$ws.Range("F40").Value = 0.25
$ws.Range("F44").Value = 0.5

# Move the active selection/view to F41 (also updates the frozen pane's
# topLeftCell to scroll the view down).
$ws.Range("F41").Select()
